$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 12785.286
$ws.Range("I18").Value = 11749.5
$ws.Range("K18").Value = 11749.5
$ws.Range("M18").Value = -11465.5

$ws.Range("H53").Value = 297.33334
$ws.Range("I53").Value = 196.57143
$ws.Range("K53").Value = 196.57143
$ws.Range("M53").Value = 440.42857

$ws.Range("H121").Value = 1714.762
$ws.Range("J121").Value = 1714.762
$ws.Range("L121").Value = 5144.286
$ws.Range("N121").Value = -8638.286

$ws.Range("H132").Value = 2568584
$ws.Range("I132").Value = 4905.9653
$ws.Range("J132").Value = 10003250
$ws.Range("K132").Value = 14717.8959
$ws.Range("L132").Value = 30009750
$ws.Range("M132").Value = -12187.8959
$ws.Range("N132").Value = -30014810

$ws.Range("H137").Value = 752795
$ws.Range("I137").Value = 1216727.2
$ws.Range("J137").Value = 3366
$ws.Range("K137").Value = 3650181.6
$ws.Range("L137").Value = 10098
$ws.Range("M137").Value = -3647631.6
$ws.Range("N137").Value = -15198

$ws.Range("H138").Value = 159258.1
$ws.Range("I138").Value = 628856.7
$ws.Range("J138").Value = 5171.0625
$ws.Range("K138").Value = 1886570.1
$ws.Range("L138").Value = 15513.1875
$ws.Range("M138").Value = -1881430.1
$ws.Range("N138").Value = -25793.1875

$ws.Range("H141").Value = 4237.2905
$ws.Range("I141").Value = 3667.4482
$ws.Range("K141").Value = 11002.3446
$ws.Range("M141").Value = -5822.3446

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H18").Value = 24994
$ws.Range("J18").Value = 24994
$ws.Range("L18").Value = 24994
$ws.Range("N18").Value = -25638

$ws.Range("H30").Value = 4566.8335
$ws.Range("I30").Value = 525.25
$ws.Range("K30").Value = 525.25
$ws.Range("M30").Value = -375.25

$ws.Range("H32").Value = 20399
$ws.Range("I32").Value = 18152.2
$ws.Range("J32").Value = 33879.8
$ws.Range("K32").Value = 18152.2
$ws.Range("L32").Value = 33879.8
$ws.Range("M32").Value = -17865.2
$ws.Range("N32").Value = -34453.8

$ws.Range("H92").Value = 275038750
$ws.Range("J92").Value = 275038750
$ws.Range("L92").Value = 275038750
$ws.Range("N92").Value = -275043742

$ws.Range("H102").Value = 6927.706
$ws.Range("I102").Value = 6718.7896
$ws.Range("K102").Value = 6718.7896
$ws.Range("M102").Value = -5096.7896

$ws.Range("H110").Value = 2584.111
$ws.Range("I110").Value = 2137.5386
$ws.Range("K110").Value = 2137.5386
$ws.Range("M110").Value = -92.53859999999986

$ws.Range("H122").Value = 2734113.2
$ws.Range("I122").Value = 6281
$ws.Range("K122").Value = 18843
$ws.Range("M122").Value = -16393

$ws.Range("H132").Value = 3233.125
$ws.Range("I132").Value = 2727.5
$ws.Range("K132").Value = 8182.5
$ws.Range("M132").Value = -5652.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 268.36365
$ws.Range("I22").Value = 243.75
$ws.Range("K22").Value = 243.75
$ws.Range("M22").Value = -70.75

$ws.Range("H94").Value = 629.37036
$ws.Range("I94").Value = 655.76
$ws.Range("K94").Value = 655.76
$ws.Range("M94").Value = -204.76

$ws.Range("H134").Value = 2500.75
$ws.Range("I134").Value = 1601.8695
$ws.Range("K134").Value = 4805.6085
$ws.Range("M134").Value = -2270.6085

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 3000000
$ws.Range("J4").Value = 3000000
$ws.Range("L4").Value = 3000000
$ws.Range("N4").Value = -3000224

$ws.Range("H22").Value = 688.8
$ws.Range("I22").Value = 355.2
$ws.Range("J22").Value = 1022.4
$ws.Range("K22").Value = 355.2
$ws.Range("L22").Value = 1022.4
$ws.Range("M22").Value = -5.199999999999989
$ws.Range("N22").Value = -1722.4

$ws.Range("H31").Value = 2267.07
$ws.Range("I31").Value = 1560.4783
$ws.Range("J31").Value = 2745.0588
$ws.Range("K31").Value = 1560.4783
$ws.Range("L31").Value = 2745.0588
$ws.Range("M31").Value = -1265.4783
$ws.Range("N31").Value = -3335.0588

$ws.Range("H34").Value = 2267.07
$ws.Range("I34").Value = 1560.4783
$ws.Range("J34").Value = 2745.0588
$ws.Range("K34").Value = 1560.4783
$ws.Range("L34").Value = 2745.0588
$ws.Range("M34").Value = -1358.4783
$ws.Range("N34").Value = -3149.0588

$ws.Range("H141").Value = 412524.62
$ws.Range("J141").Value = 545559.4
$ws.Range("L141").Value = 545559.4
$ws.Range("N141").Value = -555919.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 43144136
$ws.Range("I4").Value = 40538284
$ws.Range("K4").Value = 121614852
$ws.Range("M4").Value = -121614740

$ws.Range("H107").Value = 1800.6111
$ws.Range("I107").Value = 851
$ws.Range("J107").Value = 2165.8462
$ws.Range("K107").Value = 2553
$ws.Range("L107").Value = 6497.5386
$ws.Range("M107").Value = -633
$ws.Range("N107").Value = -10337.5386

$ws.Range("H129").Value = 1598.7778
$ws.Range("J129").Value = 3000
$ws.Range("L129").Value = 9000
$ws.Range("N129").Value = -19000

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 1015
$ws.Range("I9").Value = 1015
$ws.Range("K9").Value = 1015
$ws.Range("M9").Value = -845

$ws.Range("H14").Value = 19119334
$ws.Range("I14").Value = 22942000
$ws.Range("J14").Value = 6000
$ws.Range("K14").Value = 22942000
$ws.Range("L14").Value = 6000
$ws.Range("M14").Value = -22941832
$ws.Range("N14").Value = -6336

$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").ClearContents()

$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("M70").ClearContents()
$ws.Range("N70").ClearContents()

$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("M73").ClearContents()
$ws.Range("N73").ClearContents()

$ws.Range("H80").Value = 15326
$ws.Range("I80").Value = 20909.25
$ws.Range("J80").Value = 7881.6665
$ws.Range("K80").Value = 20909.25
$ws.Range("L80").Value = 7881.6665
$ws.Range("M80").Value = -19911.25
$ws.Range("N80").Value = -9877.666499999999

$ws.Range("H83").Value = 15326
$ws.Range("I83").Value = 20909.25
$ws.Range("J83").Value = 7881.6665
$ws.Range("K83").Value = 104546.25
$ws.Range("L83").Value = 39408.3325
$ws.Range("M83").Value = -99554.25
$ws.Range("N83").Value = -49392.3325

$ws.Range("H123").Value = 45000
$ws.Range("J123").Value = 45000
$ws.Range("L123").Value = 45000
$ws.Range("N123").Value = -49900

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 9473.368
$ws.Range("I16").Value = 9887.941000000001
$ws.Range("J16").Value = 5949.5
$ws.Range("K16").Value = 9887.941000000001
$ws.Range("L16").Value = 5949.5
$ws.Range("M16").Value = -9717.941000000001
$ws.Range("N16").Value = -6289.5

$ws.Range("H100").Value = 3228.7058
$ws.Range("I100").Value = 1090
$ws.Range("K100").Value = 1090
$ws.Range("M100").Value = -549

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 52899.75
$ws.Range("I3").Value = 70299.664
$ws.Range("J3").Value = 700
$ws.Range("K3").Value = 70299.664
$ws.Range("L3").Value = 700
$ws.Range("M3").Value = -70185.664
$ws.Range("N3").Value = -928

$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("M6").ClearContents()

$ws.Range("H44").Value = 19534.5
$ws.Range("J44").Value = 24069
$ws.Range("L44").Value = 24069
$ws.Range("N44").Value = -25177

$ws.Range("H50").Value = 14500
$ws.Range("J50").Value = 14500
$ws.Range("L50").Value = 14500
$ws.Range("N50").Value = -15762

$ws.Range("H55").Value = 4499.5
$ws.Range("J55").Value = 4499.5
$ws.Range("L55").Value = 4499.5
$ws.Range("N55").Value = -5053.5

$ws.Range("H59").Value = 24150
$ws.Range("J59").Value = 24150
$ws.Range("L59").Value = 24150
$ws.Range("N59").Value = -25626

$ws.Range("H100").Value = 34128.184
$ws.Range("J100").Value = 54116.5
$ws.Range("L100").Value = 108233
$ws.Range("N100").Value = -109315

$ws.Range("H107").Value = 8775.789000000001
$ws.Range("I107").Value = 1378.8667
$ws.Range("K107").Value = 4136.6001
$ws.Range("M107").Value = -2216.6001

$ws.Range("H122").Value = 5840.1514
$ws.Range("I122").Value = 4565.375
$ws.Range("K122").Value = 13696.125
$ws.Range("M122").Value = -11246.125

$ws.Range("H126").Value = 34556.383
$ws.Range("I126").Value = 41923.8
$ws.Range("K126").Value = 125771.4
$ws.Range("M126").Value = -123301.4
